# The document contains two "<id>...</id>" tags that were each split
# across three separate runs (e.g. "<id>" / "p089v_1" / "</id>" with
# differing formatting on the middle run). Re-downloading/re-importing
# them collapsed each trio back into a single run. We reproduce that by
# running a Find & Replace of the full tag text over itself: Word's
# find/replace merges the matched range into one run using the
# formatting of the first run in the match (Courier New / dark-yellow),
# which is exactly the target formatting for the merged run.

$d = $word.ActiveDocument

$rng = $d.Content
[void]$rng.Find.Execute("<id>p089v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p089v_1</id>", 2)

$rng2 = $d.Content
[void]$rng2.Find.Execute("<id>p089v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p089v_2</id>", 2)
